$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add a new column AQ (26-jul) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("AP1").Copy()
$ws1.Range("AQ1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("AQ1").Value = "26-jul"

$spotValues = @{
    2  = 96.27
    3  = 77.16
    4  = 66.28
    5  = 56.84
    6  = 46.4
    7  = 50.74
    8  = 57.33
    9  = 55.01
    10 = 41.05
    11 = 33.2
    12 = 17.69
    13 = 13.78
    14 = 18
    15 = 14.58
    16 = 8.67
    17 = 8.15
    18 = 5.98
    19 = 17.3
    20 = 27.13
    21 = 49.33
    22 = 65.29000000000001
    23 = 92.56999999999999
    24 = 100
    25 = 90.66
}

foreach ($row in $spotValues.Keys) {
    $ws1.Cells.Item($row, 43).Value = $spotValues[$row]
}

# --- Sheet "Gaz": append a new row 40 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2Date = $ws2.Range("A40")
$ws2Date.NumberFormat = "@"
$ws2Date.Value = "2025-07-24"
$ws2Date.ClearFormats()
$ws2.Range("B40").Value = 32.075

# --- Sheet "CO2": append a new row 40 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3Date = $ws3.Range("A40")
$ws3Date.NumberFormat = "@"
$ws3Date.Value = "2025-07-24"
$ws3Date.ClearFormats()
$ws3.Range("B40").Value = 70.2
